$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenes")

# --- Row 2 / C2: rework the "bringing back enough" phrasing into a templated
#     gerund placeholder ({industry:hometown:goodsGer}) ---
$ws.Range("C2").Value = "{name}, you live in {location:hometown:name}, a small village {location:hometown:feature:relativeposition}. You work {industry:hometown:workplace} every day, {industry:hometown:goodsGer} enough {industry:hometown:goods} {industry:hometown:purpose} your small village.`n`nOne day, a messenger rides into town. ""Citizens of {location:hometown:name}, the {character:antagonist:baron:baron} {character:baron:name} lays claim to your city. You are now under {character:baron:possPronoun} rule and will pay taxes and fealty to {character:baron:objPronoun}.""`n`nWhat do you do?"

# --- Row 27 / C27 & E27: "Chief" literal -> {chief} template placeholder
#     (so a "chieftess" variant can be substituted too) ---
$ws.Range("C27").Value = "{location:hometown:namewiththe:cap} is now free! Your fellow citizens are eagerly offering you a position of leadership: they want to make you the {chief} of {location:hometown:namewiththe}. What do you do? {|SET:chasedbybaron:false|}"
$ws.Range("E27").Value = "Take the new job as {chief} of {location:hometown:name}"

# --- Row 2 grew taller to fit the new text ---
$ws.Rows.Item(2).RowHeight = 135.8

# --- Reset the sheet view back to the top and select C2 ---
$ws.Range("C2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
